$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.359.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.218.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "109.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.58%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.07%  "

$ws.Range("E15").Value = "  -3.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.549.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.231.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.241.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("E21").Value = "  -5.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "228.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.76%  "

$ws.Range("E29").Value = "  -2.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("E31").Value = "  -6.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0865"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.125"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0364"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.18%  "

$ws.Range("E40").Value = "  -5.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.227"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.46%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.91%  "
